$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.112.13'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.239.19'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.04'
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '86.87'
$ws.Range("E6").Value = '  +4.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.85'
$ws.Range("E10").Value = '  +5.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0792'
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.99'
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.109'
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.41'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.573.60'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.13'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.198.88'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.731'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.998.22'
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.17'
$ws.Range("E21").Value = '  +9.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.81'
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.33'
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.38'
$ws.Range("E24").Value = '  +3.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.47'
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.85'
$ws.Range("E27").Value = '  +4.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.85'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("E29").Value = '  +2.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.32'
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.27'
$ws.Range("E31").Value = '  +4.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.54'
$ws.Range("E32").Value = '  +3.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0724'
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.24'
$ws.Range("E37").Value = '  +9.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.80'
$ws.Range("E38").Value = '  +4.58%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.100'
$ws.Range("E39").Value = '  +4.71%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("E41").Value = '  +5.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.041.63'
$ws.Range("E43").Value = '  +6.75%  '
$ws.Range("E44").Value = '  +7.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0271'
$ws.Range("E45").Value = '  +4.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.10'
$ws.Range("E46").Value = '  +12.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.02'
$ws.Range("E47").Value = '  +7.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.59'
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.461.94'
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.29'
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '89.44'
$ws.Range("E51").Value = '  +2.55%  '
